$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 24 (existing rows 24..31 shift down to 25..32).
#    Excel's default row-insert behaviour copies the formatting of the row above
#    (row 23, the yellow highlighted "Forderungen..." row) into the freshly
#    inserted row, which is exactly what the target workbook expects for the
#    new row 24.
$ws.Rows.Item(24).Insert()

# 2. Populate the newly inserted row 24 with the "Forderungen gegen
#    Gewaehrtraeger" entry. The numeric values used to live in (old) row 23;
#    they now belong to this new row.
$ws.Range("A24").Value = "Umlaufvermoegen"
$ws.Range("B24").Value = "Forderungen und sonstige Vermoegensgegenstaende"
$ws.Range("C24").Value = "Forderungen gegen Gewaehrtraeger"
$ws.Range("D24").Value = 795053.85
$ws.Range("E24").Value = 3360

# 3. The original row 23 keeps its label ("Forderungen gegen Unternehmen, mit
#    denen ein Beteiligungsverhaeltnis besteht") but no longer carries any
#    figures, and loses its special yellow highlighting (that now belongs to
#    the new row 24 inserted above).
$ws.Range("D23").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("A23:C23").Interior.ColorIndex = -4142
$ws.Range("A23:C23").Interior.Pattern = -4142

# 4. Move the three existing comments down by one row, to follow the data
#    that they annotate.
$oldA31 = $ws.Range("A31").Comment
$textA31 = $oldA31.Text()
$oldA31.Delete()
$newA32 = $ws.Range("A32").AddComment($textA31)

$oldB27 = $ws.Range("B27").Comment
$textB27 = $oldB27.Text()
$oldB27.Delete()
$newB28 = $ws.Range("B28").AddComment($textB27)

$oldC23 = $ws.Range("C23").Comment
$textC23 = $oldC23.Text()
$oldC23.Delete()
$newC24 = $ws.Range("C24").AddComment($textC23)

Write-Host "Done"
